# Generate Report for Handback
# Renames the in-flight handback file from "84b9a225-...-98900a48dad2" to
# "6996298e-...-04f8a2062c52" (updating its timestamps) and records a brand
# new handback file "dd5a7bb6-...-53e74a0b6ef0" as an additional row on every
# sheet (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$oldUuid = "84b9a225-6938-4a42-9132-98900a48dad2"
$uuid1   = "6996298e-972a-422f-874d-04f8a2062c52"
$uuid2   = "dd5a7bb6-83dc-49c5-b38f-53e74a0b6ef0"

$hash1 = "a98dc62e5417b615d9831d9716010847cc3c2164"
$hash2 = "022f0846ef4cd392eb68e416f92ebbffcf23b4a4"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# -- update the existing (renamed) handback file's row --
$wsOverview.Range("A2").Value = "$uuid1.md"
$wsOverview.Range("B2").Value = "e2e\$uuid1.md"
$wsOverview.Range("C2").Value = ".md"
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G2").Value = "2016-08-18 19:05:52"

# -- append the new handback file's row via the table --
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()
$rngOverview = $rowOverview.Range
$rngOverview.Cells.Item(1, 1).Value = "$uuid2.md"
$rngOverview.Cells.Item(1, 2).Value = "e2e\$uuid2.md"
$rngOverview.Cells.Item(1, 3).Value = ".md"
$rngOverview.Cells.Item(1, 5).Value = "Handed back: in sync with en-US"
$rngOverview.Cells.Item(1, 6).Value = "Handed back: in sync with en-US"
$rngOverview.Cells.Item(1, 7).Value = "2016-08-18 19:05:52"

# -- rebuild hyperlinks (B2 updated target, B3 brand new) --
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9f9632f9ccb6aff520c57a53b776ecd0d0e4429/e2e/$uuid1.md", $null, $null, "e2e\$uuid1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9f9632f9ccb6aff520c57a53b776ecd0d0e4429/e2e/$uuid2.md", $null, $null, "e2e\$uuid2.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# -- update the existing (renamed) handback file's row --
$wsZhCn.Range("A2").Value = "$uuid1.md"
$wsZhCn.Range("B2").Value = ".md"
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D2").Value = "e2e"
$wsZhCn.Range("E2").Value = "ht"
$wsZhCn.Range("F2").Value = "False"
$wsZhCn.Range("G2").Value = "$uuid1.$hash1.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-18 19:05:47"
$wsZhCn.Range("I2").Value = "$uuid1.md"
$wsZhCn.Range("J2").Value = "$uuid1.$hash1.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-18 19:06:21"
$wsZhCn.Range("L2").Value = ""
$wsZhCn.Range("M2").Value = "True"
$wsZhCn.Range("N2").Value = ""
$wsZhCn.Range("O2").Value = "False"
$wsZhCn.Range("P2").Value = ""

# -- append the new handback file's row via the table --
$loZhCn = $wsZhCn.ListObjects.Item(1)
$rowZhCn = $loZhCn.ListRows.Add()
$rngZhCn = $rowZhCn.Range
$rngZhCn.Cells.Item(1, 1).Value  = "$uuid2.md"
$rngZhCn.Cells.Item(1, 2).Value  = ".md"
$rngZhCn.Cells.Item(1, 3).Value  = "Handed back: in sync with en-US"
$rngZhCn.Cells.Item(1, 4).Value  = "e2e"
$rngZhCn.Cells.Item(1, 5).Value  = "ht"
$rngZhCn.Cells.Item(1, 6).Value  = "True"
$rngZhCn.Cells.Item(1, 7).Value  = "$uuid2.$hash2.zh-cn.xlf"
$rngZhCn.Cells.Item(1, 8).Value  = "2016-08-18 19:05:47"
$rngZhCn.Cells.Item(1, 9).Value  = "$uuid2.md"
$rngZhCn.Cells.Item(1, 10).Value = "$uuid2.$hash2.zh-cn.xlf"
$rngZhCn.Cells.Item(1, 11).Value = "2016-08-18 19:06:21"
$rngZhCn.Cells.Item(1, 12).Value = ""
$rngZhCn.Cells.Item(1, 13).Value = "True"
$rngZhCn.Cells.Item(1, 14).Value = ""
$rngZhCn.Cells.Item(1, 15).Value = "False"
$rngZhCn.Cells.Item(1, 16).Value = ""

# -- rebuild hyperlinks (A2/I2 updated target, A3/I3 brand new) --
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9f9632f9ccb6aff520c57a53b776ecd0d0e4429/e2e/$uuid1.md", $null, $null, "$uuid1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0cfe3722b9cfb02c88474925f0e8bbb26030ac10/e2e/$uuid1.md", $null, $null, "$uuid1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9f9632f9ccb6aff520c57a53b776ecd0d0e4429/e2e/$uuid2.md", $null, $null, "$uuid2.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0cfe3722b9cfb02c88474925f0e8bbb26030ac10/e2e/$uuid2.md", $null, $null, "$uuid2.md")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# -- update the existing (renamed) handback file's row --
$wsDeDe.Range("A2").Value = "$uuid1.md"
$wsDeDe.Range("B2").Value = ".md"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D2").Value = "e2e"
$wsDeDe.Range("E2").Value = "ht"
$wsDeDe.Range("F2").Value = "False"
$wsDeDe.Range("G2").Value = "$uuid1.$hash1.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-18 19:05:52"
$wsDeDe.Range("I2").Value = "$uuid1.md"
$wsDeDe.Range("J2").Value = "$uuid1.$hash1.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-18 19:06:30"
$wsDeDe.Range("L2").Value = ""
$wsDeDe.Range("M2").Value = "True"
$wsDeDe.Range("N2").Value = ""
$wsDeDe.Range("O2").Value = "False"
$wsDeDe.Range("P2").Value = ""

# -- append the new handback file's row via the table --
$loDeDe = $wsDeDe.ListObjects.Item(1)
$rowDeDe = $loDeDe.ListRows.Add()
$rngDeDe = $rowDeDe.Range
$rngDeDe.Cells.Item(1, 1).Value  = "$uuid2.md"
$rngDeDe.Cells.Item(1, 2).Value  = ".md"
$rngDeDe.Cells.Item(1, 3).Value  = "Handed back: in sync with en-US"
$rngDeDe.Cells.Item(1, 4).Value  = "e2e"
$rngDeDe.Cells.Item(1, 5).Value  = "ht"
$rngDeDe.Cells.Item(1, 6).Value  = "True"
$rngDeDe.Cells.Item(1, 7).Value  = "$uuid2.$hash2.de-de.xlf"
$rngDeDe.Cells.Item(1, 8).Value  = "2016-08-18 19:05:52"
$rngDeDe.Cells.Item(1, 9).Value  = "$uuid2.md"
$rngDeDe.Cells.Item(1, 10).Value = "$uuid2.$hash2.de-de.xlf"
$rngDeDe.Cells.Item(1, 11).Value = "2016-08-18 19:06:30"
$rngDeDe.Cells.Item(1, 12).Value = ""
$rngDeDe.Cells.Item(1, 13).Value = "True"
$rngDeDe.Cells.Item(1, 14).Value = ""
$rngDeDe.Cells.Item(1, 15).Value = "False"
$rngDeDe.Cells.Item(1, 16).Value = ""

# -- rebuild hyperlinks (A2/I2 updated target, A3/I3 brand new) --
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9f9632f9ccb6aff520c57a53b776ecd0d0e4429/e2e/$uuid1.md", $null, $null, "$uuid1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/11279c7658a5541c9570cfcb188a6596b2d298c2/e2e/$uuid1.md", $null, $null, "$uuid1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9f9632f9ccb6aff520c57a53b776ecd0d0e4429/e2e/$uuid2.md", $null, $null, "$uuid2.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/11279c7658a5541c9570cfcb188a6596b2d298c2/e2e/$uuid2.md", $null, $null, "$uuid2.md")

Write-Output "Handback report rows updated."
